# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp message (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Julio de 2020 a las 10:18"

# --- Fix swapped country labels (Islas Malvinas / Groenlandia) ---
# Row 209 previously showed "Islas Malvinas" with Row 210 showing "Groenlandia";
# the names were mislabeled, so swap them while keeping the case counts in place.
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Updated case numbers (country rows) ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 3545257
$ws.Cells.Item(4, 3).Value = 180
$ws.Cells.Item(4, 5).Value = 1805791

# Row 7: Rusia
$ws.Cells.Item(7, 2).Value = 746369
$ws.Cells.Item(7, 3).Value = 6422
$ws.Cells.Item(7, 4).Value = 523249
$ws.Cells.Item(7, 5).Value = 211350
$ws.Cells.Item(7, 7).Value = 156
$ws.Cells.Item(7, 8).Value = 11770

# Row 44: Singapur
$ws.Cells.Item(44, 2).Value = 46878
$ws.Cells.Item(44, 3).Value = 249
$ws.Cells.Item(44, 5).Value = 4114

# Row 47: Polonia
$ws.Cells.Item(47, 4).Value = 28492
$ws.Cells.Item(47, 5).Value = 8377

# Row 117: Estonia
$ws.Cells.Item(117, 2).Value = 2016
$ws.Cells.Item(117, 3).Value = 1
$ws.Cells.Item(117, 4).Value = 1901
$ws.Cells.Item(117, 5).Value = 46

# Row 118: Eslovaquia
$ws.Cells.Item(118, 2).Value = 1927
$ws.Cells.Item(118, 3).Value = 19
$ws.Cells.Item(118, 4).Value = 1507
$ws.Cells.Item(118, 5).Value = 392

# Row 137: Letonia
$ws.Cells.Item(137, 2).Value = 1178
$ws.Cells.Item(137, 3).Value = 4
$ws.Cells.Item(137, 5).Value = 125

# Row 204: Nueva Caledonia
$ws.Cells.Item(204, 2).Value = 22
$ws.Cells.Item(204, 3).Value = 1
$ws.Cells.Item(204, 5).Value = 1
